# US-13314 [IMP] PO/Catalogue mismatch: Added 2 columns to the export
#   - "PO Subtotal"        inserted right after "PO Price" (which is renamed "PO Unit Price")
#   - "Catalogue Subtotal" inserted right after "Catalogue Price"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the existing "PO Price" header to "PO Unit Price" ---
$ws.Range("F10").Value = "PO Unit Price"

# --- Insert the new "PO Subtotal" column right after the (renamed) PO Unit Price column ---
$ws.Columns("G").Insert() | Out-Null
$ws.Range("G10").Value = "PO Subtotal"
$ws.Range("G11").Value = $ws.Range("F11").Value2

# --- Insert the new "Catalogue Subtotal" column right after "Catalogue SoQ" / before it, i.e. right after Catalogue Price ---
$ws.Columns("L").Insert() | Out-Null
$ws.Range("L10").Value = "Catalogue Subtotal"

# --- Leave the selection where the author left it before saving ---
$ws.Range("I13").Select() | Out-Null
